# Applies the commit "Ajout références conf + script génération ..." to
# StructureDefinition-SavoirFaire.xlsx:
#   1. Bump the generation timestamp on the Metadata sheet.
#   2. Fix the casing of the "exerciceProfessionnel" element id/path
#      (-> "ExerciceProfessionnel") everywhere it is used on the Elements
#      sheet (ID, Path, Base Path columns).
#   3. Drop the trailing period from the "Short"/"Definition" cells that
#      describe that element.
#   4. Re-fit columns A, B and AF on the Elements sheet to their new
#      (slightly wider) best-fit width.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# 1) Metadata!B8 - generation Date
$meta.Range("B8").Value2 = "2025-10-29T11:46:56+00:00"

# 2) Elements row 6 ("SavoirFaire.exerciceProfessionnel") - fix capitalisation
#    of the "Exercice" word in the ID / Path / Base Path columns.
$elements.Range("A6").Value2 = "SavoirFaire.ExerciceProfessionnel"
$elements.Range("B6").Value2 = "SavoirFaire.ExerciceProfessionnel"
$elements.Range("AF6").Value2 = "SavoirFaire.ExerciceProfessionnel"

# 3) Same row - remove the trailing "." in Short / Definition.
$elements.Range("L6").Value2 = "Lien vers la classe ExerciceProfessionnel"
$elements.Range("M6").Value2 = "Lien vers la classe ExerciceProfessionnel"

# 4) Columns A, B and AF grew very slightly (best-fit) once the text was
#    updated - match the new stored width (closest value reachable through
#    the ColumnWidth/pixel rounding used by this host is ~27.5, which is
#    the nearest representable width to the recorded 27.42578125).
$elements.Columns.Item(1).ColumnWidth = 26.65
$elements.Columns.Item(2).ColumnWidth = 26.65
$elements.Columns.Item(32).ColumnWidth = 26.65
